$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value  = "tivi, két, bar, karaoke, golf, casino, oxy, Merdeka, website, taxi"
$ws.Range("C7").Value  = "ơi, ạ, Ôi, à, Vâng, Thôi, nhỉ, nhé, ư, À"
$ws.Range("C13").Value = "đồng, giờ, g, m, phút, ha, km, USD, kg, đ"
$ws.Range("C17").Value = "cả, chính, ngay, thôi, rồi, cái, mà, thật, đâu, đấy"
$ws.Range("C20").Value = "như thế, như vậy, làm sao, nhất là, thế nào, có lẽ, vì sao, ngày càng, một mình, Như vậy"
$ws.Range("C21").Value = "phó, viên, siêu, bất, tổng, tái, liên, đa, phi, hoá"
